$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 301.9565  # H41: was 302.47827
$ws.Cells.Item(41, 9).Value = 408  # I41: was 408.45456
$ws.Cells.Item(41, 10).Value = 204.75  # J41: was 205.33333
$ws.Cells.Item(41, 11).Value = 408  # K41: was 408.45456
$ws.Cells.Item(41, 12).Value = 204.75  # L41: was 205.33333
$ws.Cells.Item(41, 13).Value = 32  # M41: was 31.54543999999999
$ws.Cells.Item(41, 14).Value = -1084.75  # N41: was -1085.33333
$ws.Cells.Item(61, 8).Value = 912.2308  # H61: was 990.0714
$ws.Cells.Item(61, 9).Value = 705.9  # I61: was 706.1
$ws.Cells.Item(61, 10).Value = 1600  # J61: was 1700
$ws.Cells.Item(61, 11).Value = 2117.7  # K61: was 2118.3
$ws.Cells.Item(61, 12).Value = 4800  # L61: was 5100
$ws.Cells.Item(61, 13).Value = -1945.7  # M61: was -1946.3
$ws.Cells.Item(61, 14).Value = -5144  # N61: was -5444
$ws.Cells.Item(80, 8).Value = 2584.5  # H80: was 2227.3
$ws.Cells.Item(80, 9).Value = 0  # I80: was 797
$ws.Cells.Item(80, 10).Value = 2584.5  # J80: was 2386.2222
$ws.Cells.Item(80, 11).Value = 0  # K80: was 2391
$ws.Cells.Item(80, 12).Value = 7753.5  # L80: was 7158.6666
$ws.Cells.Item(80, 13).ClearContents()  # M80: was -1393
$ws.Cells.Item(80, 14).Value = -9749.5  # N80: was -9154.6666
$ws.Cells.Item(83, 8).Value = 2584.5  # H83: was 2227.3
$ws.Cells.Item(83, 9).Value = 0  # I83: was 797
$ws.Cells.Item(83, 10).Value = 2584.5  # J83: was 2386.2222
$ws.Cells.Item(83, 11).Value = 0  # K83: was 7173
$ws.Cells.Item(83, 12).Value = 23260.5  # L83: was 21475.9998
$ws.Cells.Item(83, 13).ClearContents()  # M83: was -2181
$ws.Cells.Item(83, 14).Value = -33244.5  # N83: was -31459.9998
$ws.Cells.Item(125, 8).Value = 8194.414000000001  # H125: was 8129.433
$ws.Cells.Item(125, 9).Value = 630.3889  # I125: was 644
$ws.Cells.Item(125, 10).Value = 20571.908  # J125: was 19357.584
$ws.Cells.Item(125, 11).Value = 5673.5001  # K125: was 5796
$ws.Cells.Item(125, 12).Value = 185147.172  # L125: was 174218.256
$ws.Cells.Item(125, 13).Value = -3213.5001  # M125: was -3336
$ws.Cells.Item(125, 14).Value = -190067.172  # N125: was -179138.256
$ws.Cells.Item(137, 8).Value = 2230.0557  # H137: was 2049.28
$ws.Cells.Item(137, 9).Value = 1866.5  # I137: was 1714.6154
$ws.Cells.Item(137, 11).Value = 5599.5  # K137: was 5143.8462
$ws.Cells.Item(137, 13).Value = -3049.5  # M137: was -2593.8462
$ws.Cells.Item(138, 8).Value = 8480.767  # H138: was 8751.647999999999
$ws.Cells.Item(138, 9).Value = 8376.267  # I138: was 9603.385
$ws.Cells.Item(138, 10).Value = 8506.049000000001  # J138: was 8570.130999999999
$ws.Cells.Item(138, 11).Value = 25128.801  # K138: was 28810.155
$ws.Cells.Item(138, 12).Value = 25518.147  # L138: was 25710.393
$ws.Cells.Item(138, 13).Value = -19988.801  # M138: was -23670.155
$ws.Cells.Item(138, 14).Value = -35798.147  # N138: was -35990.393
$ws.Cells.Item(141, 8).Value = 2741.04  # H141: was 2958.6086
$ws.Cells.Item(141, 9).Value = 2741.04  # I141: was 2958.6086
$ws.Cells.Item(141, 11).Value = 8223.119999999999  # K141: was 8875.825800000001
$ws.Cells.Item(141, 13).Value = -3043.119999999999  # M141: was -3695.825800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2203.0557  # H32: was 2147.5134
$ws.Cells.Item(32, 9).Value = 1291.1724  # I32: was 1222.2903
$ws.Cells.Item(32, 10).Value = 5980.857  # J32: was 6927.8335
$ws.Cells.Item(32, 11).Value = 1291.1724  # K32: was 1222.2903
$ws.Cells.Item(32, 12).Value = 5980.857  # L32: was 6927.8335
$ws.Cells.Item(32, 13).Value = -1004.1724  # M32: was -935.2902999999999
$ws.Cells.Item(32, 14).Value = -6554.857  # N32: was -7501.8335
$ws.Cells.Item(61, 8).Value = 4615.6763  # H61: was 3805.0264
$ws.Cells.Item(61, 9).Value = 4094.348  # I61: was 3228.2188
$ws.Cells.Item(61, 10).Value = 5705.727  # J61: was 6881.3335
$ws.Cells.Item(61, 11).Value = 4094.348  # K61: was 3228.2188
$ws.Cells.Item(61, 12).Value = 5705.727  # L61: was 6881.3335
$ws.Cells.Item(61, 13).Value = -3882.348  # M61: was -3016.2188
$ws.Cells.Item(61, 14).Value = -6129.727  # N61: was -7305.3335
$ws.Cells.Item(88, 8).Value = 2554.875  # H88: was 2804.875
$ws.Cells.Item(88, 10).Value = 2634.2856  # J88: was 2920
$ws.Cells.Item(88, 12).Value = 2634.2856  # L88: was 2920
$ws.Cells.Item(88, 14).Value = -3446.2856  # N88: was -3732
$ws.Cells.Item(91, 8).Value = 2554.875  # H91: was 2804.875
$ws.Cells.Item(91, 10).Value = 2634.2856  # J91: was 2920
$ws.Cells.Item(91, 12).Value = 2634.2856  # L91: was 2920
$ws.Cells.Item(91, 14).Value = -5442.2856  # N91: was -5728
$ws.Cells.Item(131, 8).Value = 0  # H131: was 75715
$ws.Cells.Item(131, 10).Value = 0  # J131: was 75715
$ws.Cells.Item(131, 12).Value = 0  # L131: was 75715
$ws.Cells.Item(131, 14).ClearContents()  # N131: was -85795
$ws.Cells.Item(136, 8).Value = 4615.6763  # H136: was 3805.0264
$ws.Cells.Item(136, 9).Value = 4094.348  # I136: was 3228.2188
$ws.Cells.Item(136, 10).Value = 5705.727  # J136: was 6881.3335
$ws.Cells.Item(136, 11).Value = 12283.044  # K136: was 9684.6564
$ws.Cells.Item(136, 12).Value = 17117.181  # L136: was 20644.0005
$ws.Cells.Item(136, 13).Value = -9733.044  # M136: was -7134.6564
$ws.Cells.Item(136, 14).Value = -22217.181  # N136: was -25744.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 14711506  # H86: was 13518722
$ws.Cells.Item(86, 9).Value = 6052.609  # I86: was 5411.923
$ws.Cells.Item(86, 11).Value = 6052.609  # K86: was 5411.923
$ws.Cells.Item(86, 13).Value = -4929.609  # M86: was -4288.923
$ws.Cells.Item(89, 8).Value = 14711506  # H89: was 13518722
$ws.Cells.Item(89, 9).Value = 6052.609  # I89: was 5411.923
$ws.Cells.Item(89, 11).Value = 30263.045  # K89: was 27059.615
$ws.Cells.Item(89, 13).Value = -24647.045  # M89: was -21443.615

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2911.8408  # H31: was 2896.0889
$ws.Cells.Item(31, 9).Value = 1603.5358  # I31: was 1618.0741
$ws.Cells.Item(31, 10).Value = 5201.375  # J31: was 4813.1113
$ws.Cells.Item(31, 11).Value = 1603.5358  # K31: was 1618.0741
$ws.Cells.Item(31, 12).Value = 5201.375  # L31: was 4813.1113
$ws.Cells.Item(31, 13).Value = -1308.5358  # M31: was -1323.0741
$ws.Cells.Item(31, 14).Value = -5791.375  # N31: was -5403.1113
$ws.Cells.Item(34, 8).Value = 2911.8408  # H34: was 2896.0889
$ws.Cells.Item(34, 9).Value = 1603.5358  # I34: was 1618.0741
$ws.Cells.Item(34, 10).Value = 5201.375  # J34: was 4813.1113
$ws.Cells.Item(34, 11).Value = 1603.5358  # K34: was 1618.0741
$ws.Cells.Item(34, 12).Value = 5201.375  # L34: was 4813.1113
$ws.Cells.Item(34, 13).Value = -1401.5358  # M34: was -1416.0741
$ws.Cells.Item(34, 14).Value = -5605.375  # N34: was -5217.1113
$ws.Cells.Item(107, 8).Value = 46047  # H107: was 46042.047
$ws.Cells.Item(107, 9).Value = 63086  # I107: was 63086.875
$ws.Cells.Item(107, 10).Value = 609.6667  # J107: was 589.1667
$ws.Cells.Item(107, 11).Value = 63086  # K107: was 63086.875
$ws.Cells.Item(107, 12).Value = 609.6667  # L107: was 589.1667
$ws.Cells.Item(107, 13).Value = -61166  # M107: was -61166.875
$ws.Cells.Item(107, 14).Value = -4449.6667  # N107: was -4429.1667
$ws.Cells.Item(134, 8).Value = 3894.5  # H134: was 2566.4
$ws.Cells.Item(134, 9).Value = 3963.3076  # I134: was 2633.0435
$ws.Cells.Item(134, 10).Value = 3000  # J134: was 1800
$ws.Cells.Item(134, 11).Value = 11889.9228  # K134: was 7899.130500000001
$ws.Cells.Item(134, 12).Value = 9000  # L134: was 5400
$ws.Cells.Item(134, 13).Value = -9354.9228  # M134: was -5364.130500000001
$ws.Cells.Item(134, 14).Value = -14070  # N134: was -10470

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 1128.2858  # H51: was 1132.1666
$ws.Cells.Item(51, 10).Value = 1181  # J51: was 1200
$ws.Cells.Item(51, 12).Value = 3543  # L51: was 3600
$ws.Cells.Item(51, 14).Value = -4463  # N51: was -4520
$ws.Cells.Item(110, 8).Value = 1829.6666  # H110: was 99.5
$ws.Cells.Item(110, 10).Value = 5290  # J110: was 0
$ws.Cells.Item(110, 12).Value = 15870  # L110: was 0
$ws.Cells.Item(110, 14).Value = -24050  # N110: was None
$ws.Cells.Item(122, 8).Value = 802.8  # H122: was 816
$ws.Cells.Item(122, 9).Value = 737.1667  # I122: was 730.75
$ws.Cells.Item(122, 11).Value = 6634.5003  # K122: was 6576.75
$ws.Cells.Item(122, 13).Value = -4184.5003  # M122: was -4126.75
$ws.Cells.Item(140, 8).Value = 169683.17  # H140: was 503250
$ws.Cells.Item(140, 9).Value = 202319.8  # I140: was 1000000
$ws.Cells.Item(140, 11).Value = 606959.3999999999  # K140: was 3000000
$ws.Cells.Item(140, 13).Value = -601779.3999999999  # M140: was -2994820

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6544.3486  # H70: was 6628.524
$ws.Cells.Item(70, 10).Value = 6265.615  # J70: was 6537
$ws.Cells.Item(70, 12).Value = 6265.615  # L70: was 6537
$ws.Cells.Item(70, 14).Value = -6805.615  # N70: was -7077
$ws.Cells.Item(73, 8).Value = 6544.3486  # H73: was 6628.524
$ws.Cells.Item(73, 10).Value = 6265.615  # J73: was 6537
$ws.Cells.Item(73, 12).Value = 6265.615  # L73: was 6537
$ws.Cells.Item(73, 14).Value = -8137.615  # N73: was -8409
$ws.Cells.Item(132, 8).Value = 1645.7727  # H132: was 1400.9354
$ws.Cells.Item(132, 9).Value = 1581.2858  # I132: was 1347.6333
$ws.Cells.Item(132, 11).Value = 4743.857400000001  # K132: was 4042.8999
$ws.Cells.Item(132, 13).Value = -2213.857400000001  # M132: was -1512.8999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3263.524  # H22: was 2976.3
$ws.Cells.Item(22, 9).Value = 2554.4  # I22: was 2534.4
$ws.Cells.Item(22, 10).Value = 3908.182  # J22: was 3197.25
$ws.Cells.Item(22, 11).Value = 2554.4  # K22: was 2534.4
$ws.Cells.Item(22, 12).Value = 3908.182  # L22: was 3197.25
$ws.Cells.Item(22, 13).Value = -2259.4  # M22: was -2239.4
$ws.Cells.Item(22, 14).Value = -4498.182  # N22: was -3787.25
$ws.Cells.Item(27, 8).Value = 3263.524  # H27: was 2976.3
$ws.Cells.Item(27, 9).Value = 2554.4  # I27: was 2534.4
$ws.Cells.Item(27, 10).Value = 3908.182  # J27: was 3197.25
$ws.Cells.Item(27, 11).Value = 2554.4  # K27: was 2534.4
$ws.Cells.Item(27, 12).Value = 3908.182  # L27: was 3197.25
$ws.Cells.Item(27, 13).Value = -2447.4  # M27: was -2427.4
$ws.Cells.Item(27, 14).Value = -4122.182  # N27: was -3411.25
$ws.Cells.Item(46, 8).Value = 2455.1538  # H46: was 2440.7856
$ws.Cells.Item(46, 9).Value = 765.3333  # I46: was 883.3333
$ws.Cells.Item(46, 10).Value = 2962.1  # J46: was 2865.5454
$ws.Cells.Item(46, 11).Value = 765.3333  # K46: was 883.3333
$ws.Cells.Item(46, 12).Value = 2962.1  # L46: was 2865.5454
$ws.Cells.Item(46, 13).Value = -577.3333  # M46: was -695.3333
$ws.Cells.Item(46, 14).Value = -3338.1  # N46: was -3241.5454
$ws.Cells.Item(55, 8).Value = 1473.2632  # H55: was 1473
$ws.Cells.Item(55, 9).Value = 268  # I55: was 287.3
$ws.Cells.Item(55, 10).Value = 3130.5  # J55: was 2790.4443
$ws.Cells.Item(55, 11).Value = 268  # K55: was 287.3
$ws.Cells.Item(55, 12).Value = 3130.5  # L55: was 2790.4443
$ws.Cells.Item(55, 13).Value = -95  # M55: was -114.3
$ws.Cells.Item(55, 14).Value = -3476.5  # N55: was -3136.4443
$ws.Cells.Item(100, 8).Value = 4922.923  # H100: was 4281.3125
$ws.Cells.Item(100, 10).Value = 3324.75  # J100: was 2827.3635
$ws.Cells.Item(100, 12).Value = 3324.75  # L100: was 2827.3635
$ws.Cells.Item(100, 14).Value = -4406.75  # N100: was -3909.3635
$ws.Cells.Item(106, 8).Value = 32500  # H106: was 0
$ws.Cells.Item(106, 10).Value = 32500  # J106: was 0
$ws.Cells.Item(106, 12).Value = 32500  # L106: was 0
$ws.Cells.Item(106, 14).Value = -35024  # N106: was None
$ws.Cells.Item(132, 8).Value = 1555.875  # H132: was 1356.619
$ws.Cells.Item(132, 9).Value = 1140.16  # I132: was 1035.2941
$ws.Cells.Item(132, 10).Value = 3040.5715  # J132: was 2722.25
$ws.Cells.Item(132, 11).Value = 3420.48  # K132: was 3105.8823
$ws.Cells.Item(132, 12).Value = 9121.7145  # L132: was 8166.75
$ws.Cells.Item(132, 13).Value = -890.4800000000005  # M132: was -575.8823000000002
$ws.Cells.Item(132, 14).Value = -14181.7145  # N132: was -13226.75
$ws.Cells.Item(136, 8).Value = 13365.312  # H136: was 13094.108
$ws.Cells.Item(136, 9).Value = 1957.3334  # I136: was 1922.9032
$ws.Cells.Item(136, 11).Value = 5872.0002  # K136: was 5768.7096
$ws.Cells.Item(136, 13).Value = -3322.0002  # M136: was -3218.7096

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1246.5  # H81: was 1105.8334
$ws.Cells.Item(81, 9).Value = 1798  # I81: was 1299
$ws.Cells.Item(81, 10).Value = 1062.6666  # J81: was 1009.25
$ws.Cells.Item(81, 11).Value = 3596  # K81: was 2598
$ws.Cells.Item(81, 12).Value = 2125.3332  # L81: was 2018.5
$ws.Cells.Item(81, 13).Value = -2535  # M81: was -1537
$ws.Cells.Item(81, 14).Value = -4247.3332  # N81: was -4140.5
$ws.Cells.Item(84, 8).Value = 1246.5  # H84: was 1105.8334
$ws.Cells.Item(84, 9).Value = 1798  # I84: was 1299
$ws.Cells.Item(84, 10).Value = 1062.6666  # J84: was 1009.25
$ws.Cells.Item(84, 11).Value = 17980  # K84: was 12990
$ws.Cells.Item(84, 12).Value = 10626.666  # L84: was 10092.5
$ws.Cells.Item(84, 13).Value = -12676  # M84: was -7686
$ws.Cells.Item(84, 14).Value = -21234.666  # N84: was -20700.5
$ws.Cells.Item(105, 8).Value = 17200  # H105: was 15666.5
$ws.Cells.Item(105, 10).Value = 17200  # J105: was 15666.5
$ws.Cells.Item(105, 12).Value = 17200  # L105: was 15666.5
$ws.Cells.Item(105, 14).Value = -24188  # N105: was -22654.5
$ws.Cells.Item(122, 8).Value = 2535.3044  # H122: was 2609.682
$ws.Cells.Item(122, 9).Value = 1606.6471  # I122: was 1650.875
$ws.Cells.Item(122, 11).Value = 4819.9413  # K122: was 4952.625
$ws.Cells.Item(122, 13).Value = -2369.9413  # M122: was -2502.625
$ws.Cells.Item(132, 8).Value = 5860.628  # H132: was 6309.275
$ws.Cells.Item(132, 9).Value = 5973.8687  # I132: was 6496.3145
$ws.Cells.Item(132, 11).Value = 17921.6061  # K132: was 19488.9435
$ws.Cells.Item(132, 13).Value = -15391.6061  # M132: was -16958.9435
$ws.Cells.Item(136, 8).Value = 1621.6111  # H136: was 1654.4
$ws.Cells.Item(136, 9).Value = 1247.8  # I136: was 1248.3
$ws.Cells.Item(136, 10).Value = 3490.6667  # J136: was 4091
$ws.Cells.Item(136, 11).Value = 3743.4  # K136: was 3744.9
$ws.Cells.Item(136, 12).Value = 10472.0001  # L136: was 12273
$ws.Cells.Item(136, 13).Value = -1193.4  # M136: was -1194.9
$ws.Cells.Item(136, 14).Value = -15572.0001  # N136: was -17373
